# customer_page start + fix E column ("shop" shared string bug -> numeric qty 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The E column (quantity) was mistakenly stored as a shared-string reference
# to the unrelated text "shop" instead of a plain number. Correct it for the
# whole data range (rows 1-130) so the cells become numeric literals; this
# also makes the "shop" shared string unused so it drops out of the workbook.
$ws.Range("E1:E130").Value = 1

# Move the selection/view to F3 (start of the "customer page" area), which
# also clears the old scrolled-down viewport (topLeftCell=A103 / C131).
$ws.Range("F3").Select() | Out-Null
